$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Sheet, $Addr, $Val)
    $r = $Sheet.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.ClearFormats()
}

Set-CellText $ws 'D2' '68.625.13'
Set-CellText $ws 'E2' '  +1.21%  '

Set-CellText $ws 'D3' '3.822.30'
Set-CellText $ws 'E3' '  +0.29%  '

Set-CellText $ws 'E4' '  +0.40%  '

Set-CellText $ws 'D5' '612.25'
Set-CellText $ws 'E5' '  +1.35%  '

Set-CellText $ws 'D6' '164.75'
Set-CellText $ws 'E6' '  -0.82%  '

Set-CellText $ws 'D7' '3.821.03'
Set-CellText $ws 'E7' '  +0.36%  '

Set-CellText $ws 'E8' '  -0.11%  '

Set-CellText $ws 'E9' '  -0.03%  '

Set-CellText $ws 'D10' '0.161'
Set-CellText $ws 'E10' '  +0.47%  '

Set-CellText $ws 'E11' '  -0.15%  '

Set-CellText $ws 'D12' '6.77'
Set-CellText $ws 'E12' '  +6.97%  '

Set-CellText $ws 'D13' '0.0000248'
Set-CellText $ws 'E13' '  -0.93%  '

Set-CellText $ws 'D14' '35.40'
Set-CellText $ws 'E14' '  -1.62%  '

Set-CellText $ws 'D15' '4.462.19'
Set-CellText $ws 'E15' '  +0.22%  '

Set-CellText $ws 'D16' '3.855.92'
Set-CellText $ws 'E16' '  +1.21%  '

Set-CellText $ws 'D17' '68.607.95'
Set-CellText $ws 'E17' '  +1.14%  '

Set-CellText $ws 'D18' '18.15'
Set-CellText $ws 'E18' '  -1.26%  '

Set-CellText $ws 'E19' '  +0.72%  '

Set-CellText $ws 'E20' '  -0.12%  '

Set-CellText $ws 'D21' '463.91'
Set-CellText $ws 'E21' '  -0.04%  '

Set-CellText $ws 'D22' '9.67'
Set-CellText $ws 'E22' '  -1.80%  '

Set-CellText $ws 'E23' '  +0.21%  '

Set-CellText $ws 'E24' '  +0.91%  '

Set-CellText $ws 'D25' '83.73'
Set-CellText $ws 'E25' '  +0.41%  '

Set-CellText $ws 'D26' '12.06'
Set-CellText $ws 'E26' '  -0.83%  '

Set-CellText $ws 'E27' '  -0.05%  '

Set-CellText $ws 'E28' '  +0.08%  '

Set-CellText $ws 'D29' '10.02'
Set-CellText $ws 'E29' '  -0.12%  '

Set-CellText $ws 'D30' '3.965.96'
Set-CellText $ws 'E30' '  +0.13%  '

Set-CellText $ws 'B31' 'ImmutableX'
Set-CellText $ws 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText $ws 'D31' '2.23'
Set-CellText $ws 'E31' '  +0.43%  '

Set-CellText $ws 'B32' 'PancakeSwap'
Set-CellText $ws 'C32' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-CellText $ws 'D32' '2.63'
Set-CellText $ws 'E32' '  -5.59%  '

Set-CellText $ws 'D33' '7.28'
Set-CellText $ws 'E33' '  -1.87%  '

Set-CellText $ws 'D34' '29.00'
Set-CellText $ws 'E34' '  -1.58%  '

Set-CellText $ws 'D35' '1.00'
Set-CellText $ws 'E35' '  -0.08%  '

Set-CellText $ws 'D36' '9.09'
Set-CellText $ws 'E36' '  +0.11%  '

Set-CellText $ws 'D37' '0.102'
Set-CellText $ws 'E37' '  +1.72%  '

Set-CellText $ws 'E38' '  +6.79%  '

Set-CellText $ws 'D39' '5.91'
Set-CellText $ws 'E39' '  +1.63%  '

Set-CellText $ws 'D40' '0.982'
Set-CellText $ws 'E40' '  -1.52%  '

Set-CellText $ws 'B41' 'dogwifhat'
Set-CellText $ws 'C41' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-CellText $ws 'D41' '3.15'
Set-CellText $ws 'E41' '  -2.52%  '

Set-CellText $ws 'B42' 'FirstDigitalUSD'
Set-CellText $ws 'C42' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-CellText $ws 'D42' '0.999'
Set-CellText $ws 'E42' '  -0.02%  '

Set-CellText $ws 'D44' '154.06'
Set-CellText $ws 'E44' '  +1.57%  '

Set-CellText $ws 'E45' '  -0.56%  '

Set-CellText $ws 'B46' 'ONDO'
Set-CellText $ws 'C46' 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-CellText $ws 'D46' '1.41'
Set-CellText $ws 'E46' '  +0.95%  '

Set-CellText $ws 'B47' 'Arweave'
Set-CellText $ws 'C47' 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-CellText $ws 'D47' '42.96'
Set-CellText $ws 'E47' '  -4.12%  '

Set-CellText $ws 'D48' '46.58'
Set-CellText $ws 'E48' '  -2.40%  '

Set-CellText $ws 'E49' '  +0.40%  '

Set-CellText $ws 'E50' '  +1.28%  '

Set-CellText $ws 'D51' '379.60'
Set-CellText $ws 'E51' '  -2.86%  '
